# Generate Report for Handoff
#
# The localization status report is regenerated: the single source file's
# status moves from "In Translation" to "Ready for handoff", and the
# handoff/generation timestamps advance a few seconds to the moment the
# report was produced. Because "Ready for handoff" is wider than
# "In Translation", the Status-related columns also grow to fit the new text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# Columns E (zh-cn) and F (de-de) hold the per-language status for the row;
# column G holds the latest handoff xliff-generation timestamp.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-20 19:06:47"

# Widen the two status columns to fit "Ready for handoff" (was sized for
# "In Translation").
$overview.Cells.Item(1, 5).ColumnWidth = 16.33
$overview.Cells.Item(1, 6).ColumnWidth = 16.33

# --- zh-cn sheet ---
# Column C is Status, column H is Latest Handoff Datetime.
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-20 19:06:43"
$zhcn.Cells.Item(1, 3).ColumnWidth = 16.33

# --- de-de sheet ---
# Column C is Status, column H is Latest Handoff Datetime.
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-20 19:06:47"
$dede.Cells.Item(1, 3).ColumnWidth = 16.33
